$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The edit inserts one new data row at the top of the data block (row 140),
# pushing the existing rows 140-200 down to 141-201. The new row re-uses the
# same "Mercado/Producto" boilerplate values shared by every data row, with a
# new date/quality/volume/price observation.

$ws.Rows(140).Insert()

$ws.Range("A140").Value2 = 7
$ws.Range("B140").Value = 'Terminal Hortofrutícola Agro Chillán'
$ws.Range("C140").Value = 'Ñuble'
$ws.Range("D140").Value2 = 44609
$ws.Range("E140").Value2 = 16
$ws.Range("F140").Value = 'Fruta'
$ws.Range("G140").Value2 = 100108
$ws.Range("H140").Value = 'Tropicales y subtropicales'
$ws.Range("I140").Value2 = 100108005
$ws.Range("J140").Value = 'Piña'
$ws.Range("K140").Value = 'Caramelo'
$ws.Range("L140").Value = 'Segunda'
$ws.Range("M140").Value2 = 60
$ws.Range("N140").Value2 = 17000
$ws.Range("O140").Value2 = 18000
$ws.Range("P140").Value2 = 17500
$ws.Range("Q140").Value = '$/caja 14 unidades'
$ws.Range("R140").Value = 'Ecuador'
$ws.Range("S140").Value2 = 1250
$ws.Range("T140").Value2 = 14
